# Weekly update: insert a new week's price record for Berenjena at
# Terminal La Palmera de La Serena. The new observation is inserted at
# row 61 (pushing the existing rows 61-97 down to 62-98) so the sheet
# stays ordered the same way it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61 - shifts existing rows 61:97 down to 62:98,
# carrying their formatting (incl. the date style on column D) along.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with this week's data.
$ws.Range("A61").Value2 = 8
$ws.Range("B61").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C61").Value2 = "Coquimbo"
$ws.Range("D61").Value2 = 44510
$ws.Range("E61").Value2 = 4
$ws.Range("F61").Value2 = 100112001
$ws.Range("G61").Value2 = "Berenjena"
$ws.Range("H61").Value2 = "Sin especificar"
$ws.Range("I61").Value2 = "Primera"
$ws.Range("J61").Value2 = 520
$ws.Range("K61").Value2 = 8000
$ws.Range("L61").Value2 = 8500
$ws.Range("M61").Value2 = 8250
$ws.Range("N61").Value2 = "$/caja 60 unidades"
$ws.Range("O61").Value2 = "Región de Arica y Parinacota"
$ws.Range("P61").Value2 = 138
$ws.Range("Q61").Value2 = 60
$ws.Range("R61").Value2 = "Hortaliza"
